# working_hours.xlsx edit
# "replaced plane matrix with admittance matrix in the power net"
#
# Functionally this commit:
#  - fixes the end time of the entry on row 142 (E142: 0:13 -> 0:30, i.e. 0.54166666666666663 -> 0.5625)
#  - inserts a brand new time-tracking entry (2014-07-21, 14:30-18:00) as the new row 143
#  - keeps one blank "entry" template row right below the data (now row 144)
#  - shifts the three summary rows (sum [min], sum [h], sum [working weeks]) down by one row
#    and updates the SUM()/derived ranges to include the new row 143

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. correct the existing row 142 end-time -----------------------------
$ws.Cells.Item(142, 5).Value = 0.5625

# --- 2. make room for the new entry: push the template/summary rows down -
# Inserting a whole row above the current "sum [min]" row (144) shifts that
# row and everything below it down by one, while row 143 (the blank
# template row) is duplicated into the newly inserted row, leaving a fresh
# blank template row in place at the new row 144 and the data row 143
# available to be filled in with the new entry below.
$ws.Cells.Item(144, 1).EntireRow.Insert()

# --- 3. fill in the new entry on row 143 ----------------------------------
$ws.Cells.Item(143, 1).Value = 2014
$ws.Cells.Item(143, 2).Value = 7
$ws.Cells.Item(143, 3).Value = 21
$ws.Cells.Item(143, 4).Value = 0.60416666666666663
$ws.Cells.Item(143, 5).Value = 0.75
$ws.Cells.Item(143, 6).Formula = "=(E143-D143)*24*60"
$ws.Cells.Item(143, 7).Formula = "=F143/60"

# --- 4. fix up the summary formulas now living on rows 145-147 -----------
$ws.Cells.Item(145, 6).Formula = "=SUM(F2:F143)"
$ws.Cells.Item(146, 6).Formula = "=F145/60"
$ws.Cells.Item(147, 6).Formula = "=F146/38.5"

# --- 5. restore the selection to match the edited workbook ---------------
$ws.Range("A144").Select() | Out-Null
